$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; D='64.155.45'; E='  +2.23%  '},
    @{Row=3; D='3.361.88'; E='  +3.38%  '},
    @{Row=4; E='  +0.03%  '},
    @{Row=5; D='528.57'; E='  +2.73%  '},
    @{Row=6; D='174.67'; E='  -2.25%  '},
    @{Row=7; E='  +0.38%  '},
    @{Row=8; D='3.359.25'; E='  +3.39%  '},
    @{Row=9; E='  -0.15%  '},
    @{Row=10; E='  -0.57%  '},
    @{Row=11; D='53.52'; E='  -6.78%  '},
    @{Row=12; E='  +3.38%  '},
    @{Row=13; E='  +1.78%  '},
    @{Row=14; D='9.10'; E='  +0.41%  '},
    @{Row=15; D='3.891.87'; E='  +3.73%  '},
    @{Row=16; D='3.353.10'; E='  +3.48%  '},
    @{Row=17; E='  +1.48%  '},
    @{Row=18; D='17.58'; E='  +0.32%  '},
    @{Row=19; D='64.077.13'; E='  +2.42%  '},
    @{Row=20; D='11.23'; E='  +3.10%  '},
    @{Row=21; E='  +2.42%  '},
    @{Row=22; D='374.33'; E='  +1.48%  '},
    @{Row=23; D='11.66'; E='  +4.03%  '},
    @{Row=24; D='4.08'; E='  +7.57%  '},
    @{Row=25; D='81.45'; E='  +2.86%  '},
    @{Row=26; D='3.73'; E='  +1.80%  '},
    @{Row=27; D='6.16'; E='  +0.45%  '},
    @{Row=28; D='2.70'; E='  +3.99%  '},
    @{Row=29; D='11.33'; E='  +0.20%  '},
    @{Row=30; D='8.28'; E='  -0.16%  '},
    @{Row=31; D='28.92'; E='  +2.18%  '},
    @{Row=32; D='632.64'; E='  +0.48%  '},
    @{Row=33; D='6.47'; E='  -3.40%  '},
    @{Row=34; D='11.24'; E='  +0.86%  '},
    @{Row=35; D='0.106'; E='  +1.33%  '},
    @{Row=36; D='57.87'; E='  -0.36%  '},
    @{Row=37; D='0.999'; E='  -0.13%  '},
    @{Row=38; D='36.46'; E='  +1.02%  '},
    @{Row=39; D='0.381'; E='  -2.93%  '},
    @{Row=40; D='0.0₃0734'; E='  +12.79%  '},
    @{Row=41; D='1.00'; E='  +0.34%  '},
    @{Row=42; E='  +9.85%  '},
    @{Row=43; D='2.978.46'; E='  +1.55%  '},
    @{Row=44; E='  +1.23%  '},
    @{Row=45; D='3.00'; E='  +7.78%  '},
    @{Row=46; E='  +4.57%  '},
    @{Row=47; E='  +2.32%  '},
    @{Row=48; D='2.62'; E='  -1.17%  '},
    @{Row=49; D='3.06'; E='  +5.49%  '},
    @{Row=50; E='  +0.87%  '},
    @{Row=51; D='137.48'; E='  +6.22%  '}
)

foreach ($u in $updates) {
    if ($u.ContainsKey("D")) {
        $c = $ws.Cells.Item($u.Row, 4)
        $c.NumberFormat = "@"
        $c.Value = $u.D
        $c.ClearFormats()
    }
    if ($u.ContainsKey("E")) {
        $c = $ws.Cells.Item($u.Row, 5)
        $c.NumberFormat = "@"
        $c.Value = $u.E
        $c.ClearFormats()
    }
}
